# Update the "matrices" score column (F) with refreshed ranking values, and
# adjust the "index" (B) and "race" (G) columns where the underlying rank
# order shifted, per the refreshed rankings computation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 14.28972879330137

# Row 3
$ws.Range("F3").Value = 13.42737676672315

# Row 4
$ws.Range("F4").Value = 8.447915419544929

# Row 5
$ws.Range("F5").Value = 7.374783189707426

# Row 6
$ws.Range("F6").Value = 6.318117057861989

# Row 7
$ws.Range("F7").Value = 6.269575784930271

# Row 8
$ws.Range("B8").Value = 33
$ws.Range("F8").Value = 5.377829192037456
$ws.Range("G8").Value = "White"

# Row 9
$ws.Range("B9").Value = 32
$ws.Range("F9").Value = 5.220260424381992
$ws.Range("G9").Value = "Black or African American"

# Row 10
$ws.Range("B10").Value = 30
$ws.Range("F10").Value = 5.0041670821604

# Row 11
$ws.Range("F11").Value = 3.086411040790361

# Row 12
$ws.Range("F12").Value = 1.224224424130261

# Row 13
$ws.Range("F13").Value = 0.3847272879396543
